$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "087-253-598048.pdf"
$ws.Range("B2").Value = "HUTZ"
$ws.Range("C2").Value = "https://www.africau.edu/images/default/sample.pdf"

$ws.Range("A3").Value = "445-310-440537.pdf"
$ws.Range("B3").Value = "HUTZ"
$ws.Range("C3").Value = "https://www.africau.edu/images/default/sample.pdf"

$ws.Range("A4").Value = "920-191-218526.pdf"
$ws.Range("B4").Value = "HUTZ"
$ws.Range("C4").Value = "https://www.africau.edu/images/default/sample.pdf"
